$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" "42.356.87"
Set-TextValue "E2" "  +0.33%  "
Set-TextValue "D3" "2.277.48"
Set-TextValue "E3" "  -0.24%  "
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "309.47"
Set-TextValue "E5" "  -3.77%  "
Set-TextValue "D6" "102.96"
Set-TextValue "E6" "  +0.76%  "
Set-TextValue "E7" "  -0.20%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "0.597"
Set-TextValue "E9" "  -0.75%  "
Set-TextValue "D10" "38.42"
Set-TextValue "E10" "  -2.21%  "
Set-TextValue "D11" "0.0894"
Set-TextValue "E11" "  -0.42%  "
Set-TextValue "D12" "8.17"
Set-TextValue "E12" "  -1.29%  "
Set-TextValue "E13" "  +0.64%  "
Set-TextValue "D14" "0.967"
Set-TextValue "E14" "  +0.81%  "
Set-TextValue "D15" "14.95"
Set-TextValue "E15" "  -0.81%  "
Set-TextValue "D16" "2.620.55"
Set-TextValue "E16" "  -0.45%  "
Set-TextValue "D17" "2.274.41"
Set-TextValue "E17" "  -0.25%  "
Set-TextValue "D18" "42.289.57"
Set-TextValue "E18" "  +0.17%  "
Set-TextValue "D19" "7.20"
Set-TextValue "E19" "  -1.86%  "
Set-TextValue "D20" "0.0000104"
Set-TextValue "E20" "  -1.41%  "
Set-TextValue "E21" "  +1.29%  "
Set-TextValue "D22" "72.51"
Set-TextValue "E22" "  -0.48%  "
Set-TextValue "E23" "  -6.03%  "
Set-TextValue "E24" "  -2.49%  "
Set-TextValue "E25" "  -3.12%  "
Set-TextValue "E26" "  +0.78%  "
Set-TextValue "B27" "Toncoin"
Set-TextValue "C27" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D27" "2.42"
Set-TextValue "E27" "  +5.02%  "
Set-TextValue "B28" "Cosmos"
Set-TextValue "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D28" "10.58"
Set-TextValue "E28" "  -2.25%  "
Set-TextValue "D29" "6.78"
Set-TextValue "E29" "  +12.26%  "
Set-TextValue "D30" "22.01"
Set-TextValue "E30" "  -1.83%  "
Set-TextValue "D31" "35.60"
Set-TextValue "E31" "  -6.32%  "
Set-TextValue "D32" "163.61"
Set-TextValue "E32" "  -0.11%  "
Set-TextValue "D33" "0.0848"
Set-TextValue "E33" "  -2.58%  "
Set-TextValue "E34" "  -2.58%  "
Set-TextValue "E35" "  +1.27%  "
Set-TextValue "E36" "  -3.79%  "
Set-TextValue "E37" "  -3.05%  "
Set-TextValue "E38" "  -2.93%  "
Set-TextValue "D39" "3.64"
Set-TextValue "E39" "  -0.66%  "
Set-TextValue "E40" "  -1.41%  "
Set-TextValue "E41" "  +0.98%  "
Set-TextValue "D42" "98.80"
Set-TextValue "E42" "  +9.43%  "
Set-TextValue "E43" "  -0.14%  "
Set-TextValue "D44" "67.99"
Set-TextValue "E44" "  -0.50%  "
Set-TextValue "D45" "0.224"
Set-TextValue "E45" "  -0.25%  "
Set-TextValue "B46" "Celestia"
Set-TextValue "C46" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D46" "11.84"
Set-TextValue "E46" "  -2.20%  "
Set-TextValue "B47" "Maker"
Set-TextValue "C47" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D47" "1.705.57"
Set-TextValue "E47" "  +6.93%  "
Set-TextValue "D48" "109.30"
Set-TextValue "E48" "  -3.53%  "
Set-TextValue "D49" "75.63"
Set-TextValue "E49" "  -5.05%  "
Set-TextValue "D50" "8.57"
Set-TextValue "E50" "  -3.96%  "
Set-TextValue "D51" "5.09"
Set-TextValue "E51" "  -2.34%  "
